$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Update existing cell values
$ws.Range("B3").Value = "SWG"
$ws.Range("B4").Value = "SWG"
$ws.Range("B10").Value = "internal/Imaging"

# 2) Insert two new rows before old row 19 (pushes old 19.. down to 21..)
$ws.Rows("19:20").Insert()

# populate new row 19
$ws.Range("A19").Value = 2
$ws.Range("B19").Value = "internal"
$ws.Range("E19").Value = "med"
$ws.Range("F19").Value = 10
$ws.Range("H19").Value = "Update Introduce authorization components (gridgrouper,  csm, etc) to new authorization support"

# populate new row 20
$ws.Range("A20").Value = 2
$ws.Range("B20").Value = "internal"
$ws.Range("E20").Value = "med"
$ws.Range("F20").Value = 20
$ws.Range("H20").Value = "Add authorization configuration support in Introduce to Resources (needed for stateful services, such as FQP, workflow, BDT, etc)"

# 3) Page setup changes
$ws.PageSetup.Orientation = 2   # xlLandscape
$ws.PageSetup.Zoom = 70
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1

$wb.Save()
